# Section 5: Toolbars and Menus completed
# Applies:
#  - Fill in COMPLETED (G) / TARGET DATE (H) date serials for the newly
#    finished rows (Section 10: Activex Controls on Worksheets, rows 45-195)
#  - Update the two "Compliance" video rows (151/152) with corrected minutes
#    and their "completed" shared-string label variants (label has extra
#    spacing before the run time, matching the pattern already used for
#    other completed rows, e.g. row 74).
#  - Move the frozen-pane view / active selection forward to reflect the
#    newly completed section.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Group 1: rows 45-53 -> both COMPLETED (G) and TARGET DATE (H) = 2017-02-22 (42788)
foreach ($r in 45..53) {
    $ws.Range("G$r").Value2 = 42788
    $ws.Range("H$r").Value2 = 42788
}

# --- Remaining groups: only TARGET DATE (H) gets filled in ---

# rows 55-62 -> 2017-02-22 (42788)
foreach ($r in 55..62) {
    $ws.Range("H$r").Value2 = 42788
}

# rows 63-76 -> 2017-02-24 (42790)
foreach ($r in 63..76) {
    $ws.Range("H$r").Value2 = 42790
}

# rows 78-93 -> 2017-02-25 (42791)
foreach ($r in 78..93) {
    $ws.Range("H$r").Value2 = 42791
}

# rows 94-101 -> 2017-02-26 (42792)
foreach ($r in 94..101) {
    $ws.Range("H$r").Value2 = 42792
}

# rows 103-116 -> 2017-02-26 (42792)
foreach ($r in 103..116) {
    $ws.Range("H$r").Value2 = 42792
}

# rows 118-134 -> 2017-02-27 (42793)
foreach ($r in 118..134) {
    $ws.Range("H$r").Value2 = 42793
}

# rows 136-150 -> 2017-02-28 (42794)
foreach ($r in 136..150) {
    $ws.Range("H$r").Value2 = 42794
}

# --- Rows 151 & 152: re-label + re-time the two Compliance videos, then mark
#     their TARGET DATE same as the rest of this group (42794).
#     NOTE: B152 must be written before B151 so the two new shared strings
#     land in the same order as the canonical workbook.
$ws.Range("B152").Value2 = "131. Compliance check part 2   9:50"
$ws.Range("C152").Value2 = 10
$ws.Range("B151").Value2 = "130. Compliance Checker for your Forms or Quizzes Part 1   3:36"
$ws.Range("C151").Value2 = 4
$ws.Range("H151").Value2 = 42794
$ws.Range("H152").Value2 = 42794

# rows 153-163 -> 2017-03-01 (42795)
foreach ($r in 153..163) {
    $ws.Range("H$r").Value2 = 42795
}

# rows 164-173 -> 2017-03-02 (42796)
foreach ($r in 164..173) {
    $ws.Range("H$r").Value2 = 42796
}

# rows 175-184 -> 2017-03-03 (42797)
foreach ($r in 175..184) {
    $ws.Range("H$r").Value2 = 42797
}

# rows 185-195 -> 2017-03-04 (42798)
foreach ($r in 185..195) {
    $ws.Range("H$r").Value2 = 42798
}

# --- Update the frozen pane / active selection to reflect scrolling down
#     into the (now complete) section.
$ws.Activate()
$ws.Application.ActiveWindow.SplitRow = 3
$ws.Range("A35").Select()
$ws.Application.ActiveWindow.FreezePanes = $true
$ws.Range("A53").Select()
